$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Semestre ideal:" value changes from "EQD-5,EQN-6" to "EQD-8,EQN-9"
#    (row 9, columns B and C both mirror the same text)
$ws.Cells.Item(9, 2).Value2 = "EQD-8,EQN-9"
$ws.Cells.Item(9, 3).Value2 = "EQD-8,EQN-9"

# 2) Add a new "Requisitos" entry: "LOQ4002 -  Reatores Quimicos  (Requisito fraco)"
#    This becomes the new row 24, pushing the existing "LOQ4055 ..." requirement
#    down to row 25. Duplicate row 24 (which already carries the correct styles
#    s="2"/s="3" and a 30pt custom row height) into a new row 25 so formatting is
#    preserved exactly, then overwrite the text of the original row 24.
$ws.Rows.Item(24).Copy()
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(25).RowHeight = 30

$newRequisito = "LOQ4002 -  Reatores Quimicos  (Requisito fraco)" + [char]10
$ws.Cells.Item(24, 2).Value2 = $newRequisito
$ws.Cells.Item(24, 3).Value2 = $newRequisito

$excel.CutCopyMode = 0
